# 2017-01-31 update: energy.gov - chunk 7
# Adds the November 2016 monthly row to Table 2.2.A, refreshes the
# "Year to Date" and "Rolling 12 Months" summary blocks with the new
# totals, and updates the title / rolling-window caption text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update title text (subtitle) ---------------------------------------
$ws.Range("A2").Value2 = "by Sector, 2006-November 2016 (Thousand Barrels)"

# --- Insert a new row for the November 2016 monthly figures -------------
# (Row 53 currently holds the "Year to Date" section header; everything
# from there down shifts down by one row.)
$ws.Rows("53:53").Insert()

# Copy the formatting from the October 2016 row (row 52) into the new
# November row (row 53) so the styles match the other monthly data rows.
$ws.Range("A52:F52").Copy()
$ws.Range("A53:F53").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A53").Value2 = "November"
$ws.Range("B53").Value2 = 1560
$ws.Range("C53").Value2 = 1198
$ws.Range("D53").Value2 = 305
$ws.Range("E53").Value2 = 11
$ws.Range("F53").Value2 = 46

# --- Refresh the "Year to Date" figures (rows 55-57 after the insert) ---
$ws.Range("B55").Value2 = 29810
$ws.Range("C55").Value2 = 18385
$ws.Range("D55").Value2 = 10321
$ws.Range("E55").Value2 = 421
$ws.Range("F55").Value2 = 684

$ws.Range("B56").Value2 = 27344
$ws.Range("C56").Value2 = 17385
$ws.Range("D56").Value2 = 9119
$ws.Range("E56").Value2 = 241
$ws.Range("F56").Value2 = 599

$ws.Range("B57").Value2 = 19315
$ws.Range("C57").Value2 = 14316
$ws.Range("D57").Value2 = 4363
$ws.Range("E57").Value2 = 119
$ws.Range("F57").Value2 = 517

# --- Update "Rolling 12 Months Ending in ..." caption (row 58) ----------
$ws.Range("A58").Value2 = "Rolling 12 Months Ending in November"

# --- Refresh the "Rolling 12 Months" figures (rows 59-60) ---------------
$ws.Range("B59").Value2 = 29065
$ws.Range("C59").Value2 = 18653
$ws.Range("D59").Value2 = 9487
$ws.Range("E59").Value2 = 272
$ws.Range("F59").Value2 = 653

$ws.Range("B60").Value2 = 20895
$ws.Range("C60").Value2 = 15493
$ws.Range("D60").Value2 = 4717
$ws.Range("E60").Value2 = "NM"
$ws.Range("F60").Value2 = 559
